$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 64
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3000
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -3496
# row 67
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3000
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -4716
# row 98
$ws.Range("H98").Value = 1111.174
$ws.Range("I98").Value = 978.9048
$ws.Range("J98").Value = 2500
$ws.Range("K98").Value = 978.9048
$ws.Range("L98").Value = 2500
$ws.Range("M98").Value = 519.0952
$ws.Range("N98").Value = -5496
# row 116
$ws.Range("H116").Value = 2848.5
$ws.Range("I116").Value = 2797.5
$ws.Range("J116").Value = 3001.5
$ws.Range("K116").Value = 2797.5
$ws.Range("L116").Value = 3001.5
$ws.Range("M116").Value = 644.5
$ws.Range("N116").Value = -9885.5
# row 122
$ws.Range("H122").Value = 1111.174
$ws.Range("I122").Value = 978.9048
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 2936.7144
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -486.7143999999998
$ws.Range("N122").Value = -12400
# row 125
$ws.Range("H125").Value = 386.53845
$ws.Range("I125").Value = 466.66666
$ws.Range("J125").Value = 362.5
$ws.Range("K125").Value = 4199.99994
$ws.Range("L125").Value = 3262.5
$ws.Range("M125").Value = -1739.99994
$ws.Range("N125").Value = -8182.5
# row 129
$ws.Range("H129").Value = 1261.6571
$ws.Range("I129").Value = 328
$ws.Range("J129").Value = 1813.3636
$ws.Range("K129").Value = 984
$ws.Range("L129").Value = 5440.0908
$ws.Range("M129").Value = 4016
$ws.Range("N129").Value = -15440.0908
# row 131
$ws.Range("H131").Value = 1797.0714
$ws.Range("I131").Value = 364.2857
$ws.Range("J131").Value = 3229.8572
$ws.Range("K131").Value = 1092.8571
$ws.Range("L131").Value = 9689.5716
$ws.Range("M131").Value = 3947.1429
$ws.Range("N131").Value = -19769.5716
# row 132
$ws.Range("H132").Value = 5126.6665
$ws.Range("I132").Value = 5085.4
$ws.Range("J132").Value = 5333
$ws.Range("K132").Value = 15256.2
$ws.Range("L132").Value = 15999
$ws.Range("M132").Value = -12726.2
$ws.Range("N132").Value = -21059
# row 135
$ws.Range("H135").Value = 801
$ws.Range("I135").Value = 371.26315
$ws.Range("J135").Value = 2161.8333
$ws.Range("K135").Value = 3341.36835
$ws.Range("L135").Value = 19456.4997
$ws.Range("M135").Value = -806.3683499999997
$ws.Range("N135").Value = -24526.4997
# row 137
$ws.Range("H137").Value = 1410.8422
$ws.Range("I137").Value = 1609.7142
$ws.Range("J137").Value = 1294.8334
$ws.Range("K137").Value = 4829.142599999999
$ws.Range("L137").Value = 3884.5002
$ws.Range("M137").Value = -2279.142599999999
$ws.Range("N137").Value = -8984.5002
# row 138
$ws.Range("H138").Value = 3086.2166
$ws.Range("I138").Value = 2782.3333
$ws.Range("J138").Value = 3288.8057
$ws.Range("K138").Value = 8346.999899999999
$ws.Range("L138").Value = 9866.417099999999
$ws.Range("M138").Value = -3206.999899999999
$ws.Range("N138").Value = -20146.4171

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 61
$ws.Range("H61").Value = 10756080
$ws.Range("I61").Value = 20836224
$ws.Range("J61").Value = 3927.4
$ws.Range("K61").Value = 20836224
$ws.Range("L61").Value = 3927.4
$ws.Range("M61").Value = -20836012
$ws.Range("N61").Value = -4351.4
# row 74
$ws.Range("H74").Value = 11113641
$ws.Range("I74").Value = 1152.0435
$ws.Range("J74").Value = 22731242
$ws.Range("K74").Value = 1152.0435
$ws.Range("L74").Value = 22731242
$ws.Range("M74").Value = -278.0435
$ws.Range("N74").Value = -22732990
# row 77
$ws.Range("H77").Value = 11113641
$ws.Range("I77").Value = 1152.0435
$ws.Range("J77").Value = 22731242
$ws.Range("K77").Value = 5760.2175
$ws.Range("L77").Value = 113656210
$ws.Range("M77").Value = -1392.2175
$ws.Range("N77").Value = -113664946
# row 97
$ws.Range("H97").Value = 1263.6111
$ws.Range("I97").Value = 817
$ws.Range("J97").Value = 3496.6667
$ws.Range("K97").Value = 817
$ws.Range("L97").Value = 3496.6667
$ws.Range("M97").Value = -321
$ws.Range("N97").Value = -4488.6667
# row 132
$ws.Range("H132").Value = 1512400.8
$ws.Range("I132").Value = 3712.973
$ws.Range("J132").Value = 5499647
$ws.Range("K132").Value = 11138.919
$ws.Range("L132").Value = 16498941
$ws.Range("M132").Value = -8608.919
$ws.Range("N132").Value = -16504001
# row 136
$ws.Range("H136").Value = 10756080
$ws.Range("I136").Value = 20836224
$ws.Range("J136").Value = 3927.4
$ws.Range("K136").Value = 62508672
$ws.Range("L136").Value = 11782.2
$ws.Range("M136").Value = -62506122
$ws.Range("N136").Value = -16882.2

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 132
$ws.Range("H132").Value = 59445
$ws.Range("I132").Value = 50000
$ws.Range("J132").Value = 62593.332
$ws.Range("K132").Value = 50000
$ws.Range("L132").Value = 62593.332
$ws.Range("M132").Value = -44940
$ws.Range("N132").Value = -72713.332

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 6055.5835
$ws.Range("I31").Value = 1732.6666
$ws.Range("J31").Value = 7496.5557
$ws.Range("K31").Value = 1732.6666
$ws.Range("L31").Value = 7496.5557
$ws.Range("M31").Value = -1437.6666
$ws.Range("N31").Value = -8086.5557
# row 34
$ws.Range("H34").Value = 6055.5835
$ws.Range("I34").Value = 1732.6666
$ws.Range("J34").Value = 7496.5557
$ws.Range("K34").Value = 1732.6666
$ws.Range("L34").Value = 7496.5557
$ws.Range("M34").Value = -1530.6666
$ws.Range("N34").Value = -7900.5557
# row 94
$ws.Range("H94").Value = 1868.375
$ws.Range("I94").Value = 1612
$ws.Range("J94").Value = 1905
$ws.Range("K94").Value = 1612
$ws.Range("L94").Value = 1905
$ws.Range("M94").Value = -1161
$ws.Range("N94").Value = -2807
# row 99
$ws.Range("H99").Value = 2070.742
$ws.Range("I99").Value = 1950
$ws.Range("J99").Value = 2099.72
$ws.Range("K99").Value = 1950
$ws.Range("L99").Value = 2099.72
$ws.Range("M99").Value = -452
$ws.Range("N99").Value = -5095.719999999999
# row 122
$ws.Range("H122").Value = 1400.2858
$ws.Range("I122").Value = 1033.2
$ws.Range("J122").Value = 1675.6
$ws.Range("K122").Value = 3099.6
$ws.Range("L122").Value = 5026.799999999999
$ws.Range("M122").Value = -649.6000000000004
$ws.Range("N122").Value = -9926.8
# row 126
$ws.Range("H126").Value = 2070.742
$ws.Range("I126").Value = 1950
$ws.Range("J126").Value = 2099.72
$ws.Range("K126").Value = 5850
$ws.Range("L126").Value = 6299.16
$ws.Range("M126").Value = -3380
$ws.Range("N126").Value = -11239.16
# row 132
$ws.Range("H132").Value = 2877
$ws.Range("I132").Value = 2627.8333
$ws.Range("J132").Value = 3624.5
$ws.Range("K132").Value = 7883.499899999999
$ws.Range("L132").Value = 10873.5
$ws.Range("M132").Value = -5353.499899999999
$ws.Range("N132").Value = -15933.5
# row 134
$ws.Range("H134").Value = 10006519
$ws.Range("I134").Value = 11911141
$ws.Range("J134").Value = 7250.75
$ws.Range("K134").Value = 35733423
$ws.Range("L134").Value = 21752.25
$ws.Range("M134").Value = -35730888
$ws.Range("N134").Value = -26822.25

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Range("H3").Value = 7183
$ws.Range("I3").Value = 4060.7693
$ws.Range("J3").Value = 10305.23
$ws.Range("K3").Value = 12182.3079
$ws.Range("L3").Value = 30915.69
$ws.Range("M3").Value = -12070.3079
$ws.Range("N3").Value = -31139.69
# row 131
$ws.Range("H131").Value = 4998.8184
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 5820.0356
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 17460.1068
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -27540.1068
# row 133
$ws.Range("H133").Value = 12571
$ws.Range("I133").Value = 11752
$ws.Range("J133").Value = 12886
$ws.Range("K133").Value = 35256
$ws.Range("L133").Value = 38658
$ws.Range("M133").Value = -30196
$ws.Range("N133").Value = -48778
# row 141
$ws.Range("H141").Value = 4510.7295
$ws.Range("I141").Value = 3226.3572
$ws.Range("J141").Value = 5292.522
$ws.Range("K141").Value = 9679.0716
$ws.Range("L141").Value = 15877.566
$ws.Range("M141").Value = -4499.071599999999
$ws.Range("N141").Value = -26237.566

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 123
$ws.Range("H123").Value = 11931.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 11931.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 11931.5
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -16831.5
# row 126
$ws.Range("H126").Value = 2440.6667
$ws.Range("I126").Value = 2440.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 7322.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4852.000100000001
$ws.Range("N126").ClearContents()
# row 132
$ws.Range("H132").Value = 3442.08
$ws.Range("I132").Value = 2473.6667
$ws.Range("J132").Value = 5932.2856
$ws.Range("K132").Value = 7421.000100000001
$ws.Range("L132").Value = 17796.8568
$ws.Range("M132").Value = -4891.000100000001
$ws.Range("N132").Value = -22856.8568

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 908.2353
$ws.Range("I46").Value = 794.3333
$ws.Range("J46").Value = 970.36365
$ws.Range("K46").Value = 794.3333
$ws.Range("L46").Value = 970.36365
$ws.Range("M46").Value = -606.3333
$ws.Range("N46").Value = -1346.36365
